$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $oldText
    $found = $find.Execute()
    if ($found) {
        $r = $find.Parent
        $r.Text = $newText
    }
    return $found
}

function InsertAfter-Text($anchorText, $toInsert) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Text = $anchorText
    $found = $find.Execute()
    if ($found) {
        $r = $find.Parent
        $r.Collapse(0)
        $r.InsertAfter($toInsert)
    }
    return $found
}

# ---------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------
Replace-Text "Rejuvenating Educational Paradigms in the Digital Era" "The Allure of Chemistry: Embarking on a Journey of Matter Transformation" | Out-Null

# ---------------------------------------------------------------
# 2. Author name paragraph (merges the 5 runs into one)
# ---------------------------------------------------------------
$pAuthor = $d.Paragraphs.Item(2).Range
$pAuthor2 = $d.Range($pAuthor.Start, $pAuthor.End - 1)
$pAuthor2.Text = "Randall Kincaid"

# ---------------------------------------------------------------
# 3. Email paragraph
# ---------------------------------------------------------------
$pEmail = $d.Paragraphs.Item(3).Range
$pEmail2 = $d.Range($pEmail.Start, $pEmail.End - 1)
$pEmail2.Text = "randallkincaid158@yahoo.com"

# ---------------------------------------------------------------
# 4. Intro / body paragraph (paragraph 5) - sentence replacements
# ---------------------------------------------------------------
Replace-Text "The educational landscape has undergone a transformative shift in the era of digital technology, reshaping pedagogical approaches and redefining the roles of both educators and learners" "In the vast realm of sciences, chemistry stands out as a captivating subject that delves into the fundamental nature of matter and its transformation" | Out-Null

Replace-Text " With unprecedented access to information and communication channels, the digital revolution has sparked a global discourse on reimagining educational paradigms to match the needs of the 21st century" " It paints a vivid picture of the intricate interplay between atoms and molecules, inviting us to unravel the secrets hidden within the molecular structure of substances" | Out-Null

InsertAfter-Text "molecular structure of substances." " Chemistry is a saga of creativity and discovery, where elements combine and rearrange themselves, revealing the wonders of the chemical reactions that shape our world. As we embark on this journey, we will witness the elegance of chemical bonding, the magic of chemical reactions, and the immense impact chemistry has on various fields, touching our lives in countless ways." | Out-Null

Replace-Text "This evolution in educational methodologies calls for a renewed focus on collaborative learning environments, prioritizing the cultivation of critical thinking skills, problem-solving abilities, and creativity" "Chemistry permeates every aspect of our existence, from the air we breathe to the food we consume and the clothes we wear" | Out-Null

Replace-Text " Educators must embrace innovative teaching strategies that leverage technology effectively to enhance engagement and knowledge retention" " It empowers us to understand the intricate mechanisms of biological processes, unlocking the mysteries of life itself" | Out-Null

Replace-Text " By integrating digital tools and resources into the curriculum, traditional boundaries can be dissolved, promoting a seamless blend of online and offline learning" " Through the lens of chemistry, we delve into the fascinating world of materials, unraveling their properties and engineering new substances with remarkable characteristics" | Out-Null

InsertAfter-Text "engineering new substances with remarkable characteristics." " The principles of chemistry guide us in harnessing energy sources, from traditional fossil fuels to innovative renewable energy technologies, enabling us to power our world sustainably." | Out-Null

Replace-Text "Moreover, the digital age presents a unique opportunity to foster global collaboration and cultural exchange" "The study of chemistry empowers us with a profound understanding of the natural world and equips us with the tools to tackle pressing global challenges" | Out-Null

Replace-Text " Virtual platforms, educational apps, and online communities empower learners worldwide to transcend geographical barriers, collaborate on projects, and engage in cross-disciplinary dialogues" " Whether it be addressing climate change, developing life-saving drugs, or engineering advanced materials for cutting-edge technologies, chemistry plays a pivotal role in shaping a better future for humanity" | Out-Null

Replace-Text " Through these digital connections, education becomes a truly global endeavor, enriching the experiences of students from diverse backgrounds and perspectives" " As we delve deeper into this captivating realm, we embark on an intellectual and practical exploration that promises to enrich our lives and open doors to endless possibilities" | Out-Null

# ---------------------------------------------------------------
# 5. Summary paragraph (paragraph 7)
# ---------------------------------------------------------------
Replace-Text "The digital revolution has brought about a paramount need to revitalize educational paradigms, shifting towards collaborative, technology-infused environments that foster critical thinking and problem-solving skills" "Chemistry is a captivating subject that explores the essence of matter and its transformation" | Out-Null

Replace-Text " By embracing digital tools, educators can create engaging and interactive learning experiences, facilitating global collaboration and cultural exchange" " It unveils the intricacies of chemical bonding, the wonders of chemical reactions, and the profound impact chemistry has on various fields, including biology, materials science, and energy" | Out-Null

Replace-Text " This transformation demands a commitment to pedagogical innovation, ensuring that education remains relevant, adaptable, and inclusive in an age driven by rapid technological advancement" " By delving into the realm of chemistry, we gain a deeper understanding of the natural world and acquire the tools to address global challenges and create a better future" | Out-Null

InsertAfter-Text "create a better future." " Chemistry empowers us to appreciate the elegance of molecular interactions, fuels our creativity, and encourages us to embrace the endless possibilities that lie within the realm of matter transformation." | Out-Null

# ---------------------------------------------------------------
# 6. Append a new empty paragraph at the end of the document
# ---------------------------------------------------------------
$endPoint = $d.Range($d.Content.End, $d.Content.End)
$endPoint.Text = "`r"

# ---------------------------------------------------------------
# 7. Fix font everywhere: TimesNewToman -> Times New Roman
#    (apply per-paragraph, excluding the paragraph mark, to avoid
#     polluting pPr/rPr)
# ---------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    if ($r.End -gt $r.Start) {
        $r2 = $d.Range($r.Start, $r.End - 1)
        if ($r2.End -gt $r2.Start) {
            $r2.Font.Name = "Times New Roman"
        }
    }
}
